$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "275.54"
Set-TextValue $ws.Range("E2") "-1.34%"
Set-TextValue $ws.Range("D3") "27.32"
Set-TextValue $ws.Range("E3") "1.40%"
Set-TextValue $ws.Range("D4") "4.785"
Set-TextValue $ws.Range("E4") "-3.27%"
Set-TextValue $ws.Range("D5") "0.06340"
Set-TextValue $ws.Range("E5") "-1.16%"
Set-TextValue $ws.Range("D6") "6.940"
Set-TextValue $ws.Range("E6") "-0.99%"
Set-TextValue $ws.Range("D7") "1.349"
Set-TextValue $ws.Range("E7") "29.01%"
Set-TextValue $ws.Range("D8") "0.8780"
Set-TextValue $ws.Range("E8") "-1.08%"
Set-TextValue $ws.Range("D9") "0.1510"
Set-TextValue $ws.Range("E9") "1.09%"
Set-TextValue $ws.Range("D10") "0.05051"
Set-TextValue $ws.Range("E10") "-2.41%"
Set-TextValue $ws.Range("D11") "0.07569"
Set-TextValue $ws.Range("E11") "2.86%"
Set-TextValue $ws.Range("D12") "0.02955"
Set-TextValue $ws.Range("E12") "-5.45%"
Set-TextValue $ws.Range("D13") "0.09024"
Set-TextValue $ws.Range("E13") "-0.49%"
Set-TextValue $ws.Range("D14") "0.001564"
Set-TextValue $ws.Range("E14") "-0.47%"
Set-TextValue $ws.Range("D15") "0.0006388"
Set-TextValue $ws.Range("E15") "1.14%"
Set-TextValue $ws.Range("D16") "0.005724"
Set-TextValue $ws.Range("E16") "-5.01%"
Set-TextValue $ws.Range("E17") "-1.39%"
Set-TextValue $ws.Range("D18") "3.297"
Set-TextValue $ws.Range("E18") "-1.76%"
Set-TextValue $ws.Range("E19") "-1.07%"
Set-TextValue $ws.Range("E20") "0.10%"
Set-TextValue $ws.Range("D21") "0.1343"
Set-TextValue $ws.Range("E21") "0.76%"
Set-TextValue $ws.Range("D22") "3.901"
Set-TextValue $ws.Range("E22") "-0.77%"
Set-TextValue $ws.Range("D23") "0.04394"
Set-TextValue $ws.Range("E23") "0.61%"
Set-TextValue $ws.Range("D24") "0.001168"
Set-TextValue $ws.Range("E24") "-0.96%"
Set-TextValue $ws.Range("D25") "0.003858"
Set-TextValue $ws.Range("E25") "4.44%"
Set-TextValue $ws.Range("D26") "0.0001199"
Set-TextValue $ws.Range("E26") "-0.32%"
Set-TextValue $ws.Range("D27") "0.0001934"
Set-TextValue $ws.Range("E27") "13.93%"
Set-TextValue $ws.Range("D40") "0.04112"
Set-TextValue $ws.Range("E40") "0.53%"
Set-TextValue $ws.Range("D41") "0.006802"
Set-TextValue $ws.Range("E41") "2.04%"
Set-TextValue $ws.Range("D42") "0.1175"
Set-TextValue $ws.Range("E42") "-0.11%"
Set-TextValue $ws.Range("D43") "0.002028"
Set-TextValue $ws.Range("E43") "-14.25%"
Set-TextValue $ws.Range("E44") "-7.99%"
Set-TextValue $ws.Range("D45") "0.00005164"
Set-TextValue $ws.Range("E45") "-1.97%"
Set-TextValue $ws.Range("D46") "1.489"
Set-TextValue $ws.Range("E46") "-36.76%"
Set-TextValue $ws.Range("D47") "0.02297"
Set-TextValue $ws.Range("E47") "2.34%"
